$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit cyclically redistributes the "Fecha"/"Volumen"/"Precio"/"Unidad" data
# among rows 2, 3, 5 and 8 (all other columns such as Mercado, Región, Categoría,
# Variedad, Calidad, Origen, Clasificación stay the same for each row):
#   new row2 <- old row5
#   new row3 <- old row2
#   new row5 <- old row8
#   new row8 <- old row3

$rows = @(2, 3, 5, 8)

# Capture the "before" values for the columns that change.
# (use Value2 for reads - Value2 correctly returns the cell's data in this
# COM-interop environment whereas bare Value returns the property accessor)
$data = @{}
foreach ($r in $rows) {
    $data[$r] = @{
        D = $ws.Range("D$r").Value2
        J = $ws.Range("J$r").Value2
        K = $ws.Range("K$r").Value2
        L = $ws.Range("L$r").Value2
        M = $ws.Range("M$r").Value2
        N = $ws.Range("N$r").Value2
        P = $ws.Range("P$r").Value2
        Q = $ws.Range("Q$r").Value2
    }
}

# Mapping of destination row -> source row (where its new data comes from)
$map = @{
    2 = 5
    3 = 2
    5 = 8
    8 = 3
}

foreach ($dest in $map.Keys) {
    $src = $map[$dest]
    $vals = $data[$src]

    $ws.Range("D$dest").Value = $vals.D
    $ws.Range("J$dest").Value = $vals.J
    $ws.Range("K$dest").Value = $vals.K
    $ws.Range("L$dest").Value = $vals.L
    $ws.Range("M$dest").Value = $vals.M
    $ws.Range("N$dest").Value = $vals.N
    $ws.Range("P$dest").Value = $vals.P
    $ws.Range("Q$dest").Value = $vals.Q
}
